$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.672.63"
$ws.Range("E2").Value = "  -2.29%  "

# Row 3
$ws.Range("D3").Value = "2.015.86"
$ws.Range("E3").Value = "  -4.46%  "

# Row 4
$ws.Range("E4").Value = "  +0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.72%  "

# Row 6
$ws.Range("E6").Value = "  +0.18%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5046"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.54%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4275"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.91%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.82"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.45%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09211"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.91%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.132"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.73%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.66"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -6.10%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.161"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -6.24%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.997.24"
$ws.Range("E14").Value = "  -6.67%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.573"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.15%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.17"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.74%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.07%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001127"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06667"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.86%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.00"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.006"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.05%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.026"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.50%  "

# Row 23
$ws.Range("D23").Value = "29.699.20"
$ws.Range("E23").Value = "  -2.35%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.05"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.85%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.277"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.88%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.37"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.75%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.86"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.38%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.609"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.65%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.360"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.22%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.74"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.65%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.066"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -7.33%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.612"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -8.50%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09992"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.50%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.900"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.91%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.784"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.92%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.692"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -8.63%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02486"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.87%  "

# Row 38
$ws.Range("E38").Value = "  -1.73%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06410"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -6.34%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6607"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -6.73%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.86"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.83%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2088"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.69%  "

# Row 43
$ws.Range("E43").Value = "  +0.11%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6389"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.96%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.231"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.07%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.56"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.48%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.292"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.67%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.538"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.98%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07051"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.52%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.147"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.57%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000321"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.12%  "
